# Insert a new week's worth of Lechuga (lettuce) price rows for
# "Femacal de La Calera" right before the existing 2021-08-13 block
# (old row 1067), shifting all the rows below it down by 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new blank rows at 1067..1072 (everything currently at/after
# row 1067 moves down by 6 rows).
$ws.Range("A1067:A1072").EntireRow.Insert()

# Common values shared by every row of this market/category block.
$mercadoId = 3
$mercado = "Femacal de La Calera"
$region = "Coquimbo"
$codreg = 5
$categoriaId = 100112033
$categoria = "Lechuga"
$calidad = "Primera"
$origen = "Provincia de Quillota"
$clasificacion = "Hortaliza"
$fecha = 44509

# Row 1067: Conconina(o)
$r = 1067
$ws.Cells.Item($r,1).Value = $mercadoId
$ws.Cells.Item($r,2).Value = $mercado
$ws.Cells.Item($r,3).Value = $region
$ws.Cells.Item($r,4).Value = $fecha
$ws.Cells.Item($r,5).Value = $codreg
$ws.Cells.Item($r,6).Value = $categoriaId
$ws.Cells.Item($r,7).Value = $categoria
$ws.Cells.Item($r,8).Value = "Conconina(o)"
$ws.Cells.Item($r,9).Value = $calidad
$ws.Cells.Item($r,10).Value = 130
$ws.Cells.Item($r,11).Value = 4000
$ws.Cells.Item($r,12).Value = 4500
$ws.Cells.Item($r,13).Value = 4269
$ws.Cells.Item($r,14).Value = "$/caja 10 unidades"
$ws.Cells.Item($r,15).Value = $origen
$ws.Cells.Item($r,16).Value = 427
$ws.Cells.Item($r,17).Value = 10
$ws.Cells.Item($r,18).Value = $clasificacion

# Row 1068: Escarola
$r = 1068
$ws.Cells.Item($r,1).Value = $mercadoId
$ws.Cells.Item($r,2).Value = $mercado
$ws.Cells.Item($r,3).Value = $region
$ws.Cells.Item($r,4).Value = $fecha
$ws.Cells.Item($r,5).Value = $codreg
$ws.Cells.Item($r,6).Value = $categoriaId
$ws.Cells.Item($r,7).Value = $categoria
$ws.Cells.Item($r,8).Value = "Escarola"
$ws.Cells.Item($r,9).Value = $calidad
$ws.Cells.Item($r,10).Value = 130
$ws.Cells.Item($r,11).Value = 5000
$ws.Cells.Item($r,12).Value = 5500
$ws.Cells.Item($r,13).Value = 5231
$ws.Cells.Item($r,14).Value = "$/caja 15 unidades"
$ws.Cells.Item($r,15).Value = $origen
$ws.Cells.Item($r,16).Value = 349
$ws.Cells.Item($r,17).Value = 15
$ws.Cells.Item($r,18).Value = $clasificacion

# Row 1069: Española
$r = 1069
$ws.Cells.Item($r,1).Value = $mercadoId
$ws.Cells.Item($r,2).Value = $mercado
$ws.Cells.Item($r,3).Value = $region
$ws.Cells.Item($r,4).Value = $fecha
$ws.Cells.Item($r,5).Value = $codreg
$ws.Cells.Item($r,6).Value = $categoriaId
$ws.Cells.Item($r,7).Value = $categoria
$ws.Cells.Item($r,8).Value = "Española"
$ws.Cells.Item($r,9).Value = $calidad
$ws.Cells.Item($r,10).Value = 105
$ws.Cells.Item($r,11).Value = 4500
$ws.Cells.Item($r,12).Value = 5000
$ws.Cells.Item($r,13).Value = 4738
$ws.Cells.Item($r,14).Value = "$/caja 18 unidades"
$ws.Cells.Item($r,15).Value = $origen
$ws.Cells.Item($r,16).Value = 263
$ws.Cells.Item($r,17).Value = 18
$ws.Cells.Item($r,18).Value = $clasificacion

# Row 1070: Francesa morada
$r = 1070
$ws.Cells.Item($r,1).Value = $mercadoId
$ws.Cells.Item($r,2).Value = $mercado
$ws.Cells.Item($r,3).Value = $region
$ws.Cells.Item($r,4).Value = $fecha
$ws.Cells.Item($r,5).Value = $codreg
$ws.Cells.Item($r,6).Value = $categoriaId
$ws.Cells.Item($r,7).Value = $categoria
$ws.Cells.Item($r,8).Value = "Francesa morada"
$ws.Cells.Item($r,9).Value = $calidad
$ws.Cells.Item($r,10).Value = 50
$ws.Cells.Item($r,11).Value = 4500
$ws.Cells.Item($r,12).Value = 4500
$ws.Cells.Item($r,13).Value = 4500
$ws.Cells.Item($r,14).Value = "$/caja 18 unidades"
$ws.Cells.Item($r,15).Value = $origen
$ws.Cells.Item($r,16).Value = 250
$ws.Cells.Item($r,17).Value = 18
$ws.Cells.Item($r,18).Value = $clasificacion

# Row 1071: Marina
$r = 1071
$ws.Cells.Item($r,1).Value = $mercadoId
$ws.Cells.Item($r,2).Value = $mercado
$ws.Cells.Item($r,3).Value = $region
$ws.Cells.Item($r,4).Value = $fecha
$ws.Cells.Item($r,5).Value = $codreg
$ws.Cells.Item($r,6).Value = $categoriaId
$ws.Cells.Item($r,7).Value = $categoria
$ws.Cells.Item($r,8).Value = "Marina"
$ws.Cells.Item($r,9).Value = $calidad
$ws.Cells.Item($r,10).Value = 60
$ws.Cells.Item($r,11).Value = 4500
$ws.Cells.Item($r,12).Value = 4500
$ws.Cells.Item($r,13).Value = 4500
$ws.Cells.Item($r,14).Value = "$/caja 18 unidades"
$ws.Cells.Item($r,15).Value = $origen
$ws.Cells.Item($r,16).Value = 250
$ws.Cells.Item($r,17).Value = 18
$ws.Cells.Item($r,18).Value = $clasificacion

# Row 1072: Milanesa
$r = 1072
$ws.Cells.Item($r,1).Value = $mercadoId
$ws.Cells.Item($r,2).Value = $mercado
$ws.Cells.Item($r,3).Value = $region
$ws.Cells.Item($r,4).Value = $fecha
$ws.Cells.Item($r,5).Value = $codreg
$ws.Cells.Item($r,6).Value = $categoriaId
$ws.Cells.Item($r,7).Value = $categoria
$ws.Cells.Item($r,8).Value = "Milanesa"
$ws.Cells.Item($r,9).Value = $calidad
$ws.Cells.Item($r,10).Value = 110
$ws.Cells.Item($r,11).Value = 4000
$ws.Cells.Item($r,12).Value = 4500
$ws.Cells.Item($r,13).Value = 4273
$ws.Cells.Item($r,14).Value = "$/caja 20 unidades"
$ws.Cells.Item($r,15).Value = $origen
$ws.Cells.Item($r,16).Value = 214
$ws.Cells.Item($r,17).Value = 20
$ws.Cells.Item($r,18).Value = $clasificacion
